# Update benchmark results in readme
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "1T" (single-thread) benchmark numbers with newer results.
$ws.Range("C3").Value2 = 100000
$ws.Range("G3").Value2 = 385.6
$ws.Range("G4").Value2 = 168.7
$ws.Range("G5").Value2 = 940.9
$ws.Range("G6").Value2 = 1728

# Rename the "per sec total 1T" header to simply "per sec" ...
$ws.Range("I2").Value2 = "per sec"

# ... and drop the "16T" column/header entirely, along with the old
# "1T"/"16T" benchmark rows (7-10) that are no longer relevant.
$ws.Range("J2").Clear()
$ws.Range("A7:A10").EntireRow.Delete()

$wb.Application.Calculate()

# Match the author's final selection/cursor position.
$ws.Range("F7").Select() | Out-Null
